$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new labels for B1:K1 (resource categories instead of the old R1/R2/... placeholders) ---
$ws.Range("B1").Value = "population"
$ws.Range("C1").Value = "metalElements"
$ws.Range("D1").Value = "timber"
$ws.Range("E1").Value = "landArea"
$ws.Range("F1").Value = "metalAlloys"
$ws.Range("G1").Value = "electronics"
$ws.Range("H1").Value = "housing"
$ws.Range("I1").Value = "metalAlloysWaste"
$ws.Range("J1").Value = "housingWaste"
$ws.Range("K1").Value = "electronicsWaste"

# Shade the header row with a light gray fill (keeps the existing centered alignment)
$ws.Range("B1:K1").Interior.Color = 15921906

# --- Column E: random starting land-area formula, volatile via RANDBETWEEN ---
$ws.Range("E2").Formula = "=SUM(PRODUCT(RANDBETWEEN(0,5), 5000),10000)"
$ws.Range("E3:E6").Formula = "=SUM(PRODUCT(RANDBETWEEN(0,5), 5000),10000)"

# --- New column K: electronics waste tracking, defaults to 0 ---
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("K6").Value = 0

$ws.Range("K11").Select()
